$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for every existing data
#    row (2-307) from 45204 (2023-10-05) to 45205 (2023-10-06).
$ws.Range("C2:C307").Value = 45205

# 2. Row 307 gains an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(307).RowHeight = 15

# 3. Append the new record as row 308.
$ws.Cells.Item(308, 1).Value = "A 47893-2023"

$ws.Cells.Item(308, 2).NumberFormat = $ws.Cells.Item(307, 2).NumberFormat
$ws.Cells.Item(308, 2).Value = 45204

$ws.Cells.Item(308, 3).NumberFormat = $ws.Cells.Item(307, 3).NumberFormat
$ws.Cells.Item(308, 3).Value = 45205

$ws.Cells.Item(308, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(308, 5).Value = "SÄVSJÖ"

$ws.Cells.Item(308, 7).Value = 1.2
$ws.Cells.Item(308, 8).Value = 0
$ws.Cells.Item(308, 9).Value = 0
$ws.Cells.Item(308, 10).Value = 0
$ws.Cells.Item(308, 11).Value = 0
$ws.Cells.Item(308, 12).Value = 0
$ws.Cells.Item(308, 13).Value = 0
$ws.Cells.Item(308, 14).Value = 0
$ws.Cells.Item(308, 15).Value = 0
$ws.Cells.Item(308, 16).Value = 0
$ws.Cells.Item(308, 17).Value = 0

$ws.Cells.Item(308, 18).WrapText = $True
$ws.Cells.Item(308, 18).Value = ""
